$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E3").Value = 13.012
$ws.Range("E14").Value = 13.81940000000002
$ws.Range("E21").Value = 12.9965
$ws.Range("E23").Value = 13.9612
$ws.Range("E25").Value = 13.0544
